# Heidi and Hanna's Bier Haus review - copy edit
# - shorten the H1 + bold recap line, re-curl the apostrophe
# - rework the "What we like" / "What we don't like" bullet lists
# - rewrite the meta-description paragraph

$d = $word.ActiveDocument
$apos = [char]0x2019

function Get-ParaIndex($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.TrimEnd() -eq $text) {
            return $i
        }
    }
    return -1
}

function Set-ParaText($doc, $oldText, $newText) {
    $idx = Get-ParaIndex $doc $oldText
    if ($idx -eq -1) {
        throw "Paragraph not found: $oldText"
    }
    $p = $doc.Paragraphs.Item($idx)
    $p.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $newText, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1) Title (Heading1) and the bold recap line later in the doc both share the
#    exact same old text - update each occurrence.
# ---------------------------------------------------------------------------
$oldTitle = "Play Heidi and Hanna's Bier Haus for Free - A Fun Oktoberfest Slot"
$newTitle = "Play Heidi and Hanna" + $apos + "s Bier Haus for Free"
Set-ParaText $d $oldTitle $newTitle
Set-ParaText $d $oldTitle $newTitle

# ---------------------------------------------------------------------------
# 2) "What we like" bullet list.
#    Before: Interesting bonus rounds and features
#            Humorous and entertaining graphics
#            Suitable for players of all skill levels
#            High RTP of 96.15%
#    After:  Humorous and entertaining graphics
#            Impressive return to player (RTP) rate of 96.15%
#            Wide range of bonus rounds and special features
#            Suitable for players of all skill levels
# ---------------------------------------------------------------------------

# Drop the "Interesting bonus rounds and features" bullet entirely (whole
# paragraph, including its paragraph mark).
$dropIdx = Get-ParaIndex $d "Interesting bonus rounds and features"
$d.Paragraphs.Item($dropIdx).Range.Delete() | Out-Null

# "Suitable for players of all skill levels" -> becomes the RTP bullet.
Set-ParaText $d "Suitable for players of all skill levels" "Impressive return to player (RTP) rate of 96.15%"

# "High RTP of 96.15%" -> becomes the bonus-rounds bullet.
Set-ParaText $d "High RTP of 96.15%" "Wide range of bonus rounds and special features"

# Append a new final bullet restoring "Suitable for players of all skill
# levels" as the last item of the "What we like" list (same ListBullet
# paragraph style/formatting as its neighbours).
$wideRangeIdx = Get-ParaIndex $d "Wide range of bonus rounds and special features"
$d.Paragraphs.Item($wideRangeIdx).Range.InsertParagraphAfter() | Out-Null
$newBulletPara = $d.Paragraphs.Item($wideRangeIdx + 1)
$newBulletPara.Range.InsertBefore("Suitable for players of all skill levels") | Out-Null

# ---------------------------------------------------------------------------
# 3) "What we don't like" bullet list.
#    Before: Limited maximum bet of 180 euros
#            May not appeal to players looking for more traditional slots
#    After:  Limited number of paylines (up to 50)
#            No progressive jackpot feature
# ---------------------------------------------------------------------------
Set-ParaText $d "Limited maximum bet of 180 euros" "Limited number of paylines (up to 50)"
Set-ParaText $d "May not appeal to players looking for more traditional slots" "No progressive jackpot feature"

# ---------------------------------------------------------------------------
# 4) Meta description paragraph (italic line at the very end of the doc).
# ---------------------------------------------------------------------------
$oldMeta = "Read our review of Heidi and Hanna's Bier Haus, a fun and entertaining online slot game set in Oktoberfest. Play for free and enjoy bonus rounds and special features."
$newMeta = "Read our review of Heidi and Hanna" + $apos + "s Bier Haus slot game and play for free."
Set-ParaText $d $oldMeta $newMeta
